# Workman Brandon 2021 save_data regen: switch column G ("K") to use the
# strikeout figure (K) pulled from the regenerated box score instead of the
# old "Strike#" count, and (re)write the recomputed s_vals into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K value, keyed by the sheet row number (header is row 1, data
# starts at row 2)
$kVals = [ordered]@{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 2
    9  = 2
    10 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 2
    18 = 2
    19 = 1
    20 = 2
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 2
    28 = 1
    29 = 2
    30 = 0
    32 = 1
    34 = 1
    35 = 2
    36 = 2
    38 = 2
    39 = 2
}

foreach ($row in $kVals.Keys) {
    $ws.Cells.Item($row, 7).Value = $kVals[$row]
}
